$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 213.14285
$ws.Range("I41").Value = 232.55556
$ws.Range("K41").Value = 232.55556
$ws.Range("M41").Value = 207.44444
$ws.Range("H62").Value = 3855.6667
$ws.Range("I62").Value = 4560.8
$ws.Range("J62").Value = 3352
$ws.Range("K62").Value = 4560.8
$ws.Range("L62").Value = 3352
$ws.Range("M62").Value = -3936.8
$ws.Range("N62").Value = -4600
$ws.Range("H65").Value = 3855.6667
$ws.Range("I65").Value = 4560.8
$ws.Range("J65").Value = 3352
$ws.Range("K65").Value = 22804
$ws.Range("L65").Value = 16760
$ws.Range("M65").Value = -19684
$ws.Range("N65").Value = -23000
$ws.Range("H92").Value = 626.7895
$ws.Range("I92").Value = 609.7857
$ws.Range("J92").Value = 674.4
$ws.Range("K92").Value = 609.7857
$ws.Range("L92").Value = 674.4
$ws.Range("M92").Value = 638.2143
$ws.Range("N92").Value = -3170.4
$ws.Range("H100").Value = 26246.488
$ws.Range("I100").Value = 34412.266
$ws.Range("J100").Value = 3976.182
$ws.Range("K100").Value = 34412.266
$ws.Range("L100").Value = 3976.182
$ws.Range("M100").Value = -33871.266
$ws.Range("N100").Value = -5058.182
$ws.Range("H107").Value = 908.8889
$ws.Range("I107").Value = 541.4
$ws.Range("K107").Value = 541.4
$ws.Range("M107").Value = 1378.6
$ws.Range("H113").Value = 4674.3335
$ws.Range("I113").Value = 4387.125
$ws.Range("J113").Value = 5248.75
$ws.Range("K113").Value = 4387.125
$ws.Range("L113").Value = 5248.75
$ws.Range("M113").Value = -1133.125
$ws.Range("N113").Value = -11756.75
$ws.Range("H118").Value = 716.4666999999999
$ws.Range("I118").Value = 534.38464
$ws.Range("K118").Value = 1603.15392
$ws.Range("M118").Value = 53.84608000000003
$ws.Range("H132").Value = 1165.2703
$ws.Range("I132").Value = 1174.6857
$ws.Range("K132").Value = 3524.0571
$ws.Range("M132").Value = -994.0571
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3170.5625
$ws.Range("I2").Value = 2284.8635
$ws.Range("K2").Value = 2284.8635
$ws.Range("M2").Value = -2171.8635
$ws.Range("H32").Value = 6468.356
$ws.Range("I32").Value = 4855.741
$ws.Range("K32").Value = 4855.741
$ws.Range("M32").Value = -4568.741
$ws.Range("H61").Value = 7607.125
$ws.Range("J61").Value = 21449
$ws.Range("L61").Value = 21449
$ws.Range("N61").Value = -21873
$ws.Range("H102").Value = 3001.2144
$ws.Range("I102").Value = 3001.3076
$ws.Range("K102").Value = 3001.3076
$ws.Range("M102").Value = -1379.3076
$ws.Range("H116").Value = 3170.5625
$ws.Range("I116").Value = 2284.8635
$ws.Range("K116").Value = 2284.8635
$ws.Range("M116").Value = 9.136500000000069
$ws.Range("H132").Value = 5452.148
$ws.Range("I132").Value = 5084.923
$ws.Range("K132").Value = 15254.769
$ws.Range("M132").Value = -12724.769
$ws.Range("H136").Value = 7607.125
$ws.Range("J136").Value = 21449
$ws.Range("L136").Value = 64347
$ws.Range("N136").Value = -69447
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3170.5625
$ws.Range("I3").Value = 2284.8635
$ws.Range("K3").Value = 2284.8635
$ws.Range("M3").Value = -2170.8635
$ws.Range("H86").Value = 591963.7
$ws.Range("I86").Value = 1668657.4
$ws.Range("J86").Value = 4676.273
$ws.Range("K86").Value = 1668657.4
$ws.Range("L86").Value = 4676.273
$ws.Range("M86").Value = -1667534.4
$ws.Range("N86").Value = -6922.273
$ws.Range("H89").Value = 591963.7
$ws.Range("I89").Value = 1668657.4
$ws.Range("J89").Value = 4676.273
$ws.Range("K89").Value = 8343287
$ws.Range("L89").Value = 23381.365
$ws.Range("M89").Value = -8337671
$ws.Range("N89").Value = -34613.36500000001
$ws.Range("H99").Value = 1702.125
$ws.Range("I99").Value = 1702.125
$ws.Range("K99").Value = 1702.125
$ws.Range("M99").Value = -204.125
$ws.Range("H107").Value = 382.80768
$ws.Range("I107").Value = 394.18182
$ws.Range("K107").Value = 394.18182
$ws.Range("M107").Value = 1525.81818
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 194.8
$ws.Range("I22").Value = 216.91667
$ws.Range("J22").Value = 106.333336
$ws.Range("K22").Value = 216.91667
$ws.Range("L22").Value = 106.333336
$ws.Range("M22").Value = 133.08333
$ws.Range("N22").Value = -806.333336
$ws.Range("H31").Value = 3441.5
$ws.Range("J31").Value = 3817
$ws.Range("L31").Value = 3817
$ws.Range("N31").Value = -4407
$ws.Range("H34").Value = 3441.5
$ws.Range("J34").Value = 3817
$ws.Range("L34").Value = 3817
$ws.Range("N34").Value = -4221
$ws.Range("H122").Value = 1530.4
$ws.Range("I122").Value = 1626.125
$ws.Range("K122").Value = 4878.375
$ws.Range("M122").Value = -2428.375
$ws.Range("H134").Value = 5336.2
$ws.Range("I134").Value = 4358.75
$ws.Range("J134").Value = 9246
$ws.Range("K134").Value = 13076.25
$ws.Range("L134").Value = 27738
$ws.Range("M134").Value = -10541.25
$ws.Range("N134").Value = -32808
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 61.80645
$ws.Range("J2").Value = 99.25
$ws.Range("L2").Value = 595.5
$ws.Range("N2").Value = -821.5
$ws.Range("H29").Value = 987.8570999999999
$ws.Range("J29").Value = 984.8
$ws.Range("L29").Value = 2954.4
$ws.Range("N29").Value = -3508.4
$ws.Range("H109").Value = 1317.6666
$ws.Range("I109").Value = 1357.375
$ws.Range("J109").Value = 1000
$ws.Range("K109").Value = 4072.125
$ws.Range("L109").Value = 3000
$ws.Range("M109").Value = -3032.125
$ws.Range("N109").Value = -5080
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4034.4167
$ws.Range("I80").Value = 3427.375
$ws.Range("J80").Value = 5248.5
$ws.Range("K80").Value = 3427.375
$ws.Range("L80").Value = 5248.5
$ws.Range("M80").Value = -2429.375
$ws.Range("N80").Value = -7244.5
$ws.Range("H83").Value = 4034.4167
$ws.Range("I83").Value = 3427.375
$ws.Range("J83").Value = 5248.5
$ws.Range("K83").Value = 17136.875
$ws.Range("L83").Value = 26242.5
$ws.Range("M83").Value = -12144.875
$ws.Range("N83").Value = -36226.5
$ws.Range("H107").Value = 259.58066
$ws.Range("I107").Value = 287.42856
$ws.Range("J107").Value = 201.1
$ws.Range("K107").Value = 287.42856
$ws.Range("L107").Value = 201.1
$ws.Range("M107").Value = 1632.57144
$ws.Range("N107").Value = -4041.1
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 80348.484
$ws.Range("I20").Value = 16005
$ws.Range("K20").Value = 16005
$ws.Range("M20").Value = -15779
$ws.Range("H22").Value = 1888
$ws.Range("I22").Value = 883.3333
$ws.Range("K22").Value = 883.3333
$ws.Range("M22").Value = -588.3333
$ws.Range("H27").Value = 1888
$ws.Range("I27").Value = 883.3333
$ws.Range("K27").Value = 883.3333
$ws.Range("M27").Value = -776.3333
$ws.Range("H46").Value = 1573.871
$ws.Range("I46").Value = 1115.3077
$ws.Range("J46").Value = 1905.0555
$ws.Range("K46").Value = 1115.3077
$ws.Range("L46").Value = 1905.0555
$ws.Range("M46").Value = -927.3077000000001
$ws.Range("N46").Value = -2281.0555
$ws.Range("H61").Value = 3399.4
$ws.Range("J61").Value = 3619.25
$ws.Range("L61").Value = 3619.25
$ws.Range("N61").Value = -4023.25
$ws.Range("H113").Value = 3399.4
$ws.Range("J113").Value = 3619.25
$ws.Range("L113").Value = 3619.25
$ws.Range("N113").Value = -7959.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 66999.336
$ws.Range("I64").Value = 65499.5
$ws.Range("K64").Value = 65499.5
$ws.Range("M64").Value = -65251.5
$ws.Range("H67").Value = 66999.336
$ws.Range("I67").Value = 65499.5
$ws.Range("K67").Value = 65499.5
$ws.Range("M67").Value = -64641.5
